# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" table (rows 16-53) lists one row per monthly period
# (column E, text codes like "1702".."2003") together with its "Valor Mora"
# (column F). The database this report is generated from was refreshed: the
# same 38 periods are still listed, but now newest-first instead of
# oldest-first, and the "Valor Mora" figures that went along with the two
# halves of the table were swapped to match (first half now 31249, second
# half now 27578).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Period codes currently in E16:E53, oldest -> newest.
$periods = @( `
    "1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712", `
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
    "2001","2002","2003")

# New order: newest -> oldest.
$newPeriods = @($periods[($periods.Length - 1)..0])

$firstRow = 16
$lastRow = 53
$halfCount = 19

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = $firstRow + $i

    # Column E: "Periodo Mora" text code for this row.
    $ws.Cells.Item($row, 5).Value = $newPeriods[$i]

    # Column F: "Valor Mora" - first 19 rows take the value the second half
    # used to have, the remaining 19 rows take the value the first half
    # used to have.
    if ($i -lt $halfCount) {
        $ws.Cells.Item($row, 6).Value = 31249
    } else {
        $ws.Cells.Item($row, 6).Value = 27578
    }
}
